$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the two new "Robustness" rows.
#    Row 22 becomes a brand-new row (pushing the old row22.."row30"/"row32"
#    block down by one); a second insert a little further down makes room
#    for the second "Robustness" row. Excel auto-shifts every formula
#    reference (G3, and the bottom total row) as part of the native
#    Insert() operation, exactly like using the Excel UI.
# ------------------------------------------------------------------
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(26).Insert()

# ------------------------------------------------------------------
# 2. Label the two new rows (same sub-section label used elsewhere:
#    "Robustness", a brand new shared string).
# ------------------------------------------------------------------
$ws.Range("C22").Value = "Robustness"
$ws.Range("C26").Value = "Robustness"

# ------------------------------------------------------------------
# 3. Add the new "M" column (new review-round date 6/18/2014) mirroring
#    the number formatting already used by the "L" column (date style on
#    row 1, plain numbers below).
# ------------------------------------------------------------------
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "6/18/2014"

$ws.Range("M2").Value = "Pages"

$ws.Range("M3").Value = 31.1
$ws.Range("M4").Value = 1.5
$ws.Range("M5").Value = 3.5
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 0.7
$ws.Range("M9").Value = 0.9
$ws.Range("M10").Value = 0.9
$ws.Range("M11").Value = 0.3
$ws.Range("M13").Value = 5.4
$ws.Range("M14").Value = 1
$ws.Range("M15").Value = 0.75
$ws.Range("M16").Value = 0.6
$ws.Range("M17").Value = 1.5
$ws.Range("M18").Value = 0.6
$ws.Range("M19").Value = 0.7
$ws.Range("M20").Value = 2.9
$ws.Range("M21").Value = 0.7
$ws.Range("M22").Value = 0.7
$ws.Range("M23").Value = 1.4
$ws.Range("M24").Value = 2.3
$ws.Range("M25").Value = 0.6
$ws.Range("M26").Value = 0.7
$ws.Range("M27").Value = 1
$ws.Range("M28").Value = 4.4
$ws.Range("M29").Value = 0.15
$ws.Range("M30").Value = 4.8
$ws.Range("M31").Value = 0
$ws.Range("M32").Value = 3

# New cross-check subtotal next to the new column (sums the "Model
# extensions" sub-rows, M14:M19).
$ws.Range("N17").Formula = "=SUM(M14:M19)"

# Grand total for the new review round, same pattern as the existing L-column total.
$ws.Range("M34").Formula = "=M4+M5+M6+M13+M20+M24+M28+M29+M30+M31+M32"

# ------------------------------------------------------------------
# 4. Match column widths: M should look like the L/I "best fit" date columns.
# ------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(12).ColumnWidth

# ------------------------------------------------------------------
# 5. Restore the cursor/selection to where the reviewer left it.
# ------------------------------------------------------------------
$ws.Range("R12").Select()
